$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert the recurring "notes" slide at the very front of the deck.
#    It is a duplicate of the existing "If you could, what would you change
#    about our project?" slide (so it keeps the exact same placeholder /
#    formatting), moved to become the new slide 1.
# ---------------------------------------------------------------------------
$noteSlide = $p.Slides.Item(2)
$dup = $noteSlide.Duplicate()
$dup.Item(1).MoveTo(1)

# After the steps above the slide order is:
#   1: If you could... (new duplicate)
#   2: Design Patterns / Flyweight   (original slide 1)
#   3: If you could...               (original slide 2)
#   4: Motivation
#   5: Parts of the Pattern
#   6: Consequences
# Move the original "If you could..." slide so the two notes slides sit
# together at the front, pushing "Design Patterns" back to slide 3.
$p.Slides.Item(3).MoveTo(2)

# ---------------------------------------------------------------------------
# 2. Update the text of the (now second) notes slide with the new prompt.
# ---------------------------------------------------------------------------
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Text = "If anyone is interested in leading a book club session, let me know"

# ---------------------------------------------------------------------------
# 3. Small copy fix on the Motivation slide: "Intrinsic:state" -> "Intrinsic: state"
# ---------------------------------------------------------------------------
$motivation = $p.Slides.Item(4)
$contentShape = $motivation.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$introPara = $tr.Paragraphs(4, 1)
$len = $introPara.Text.Length
$introPara.Characters(1, $len).Text = "Intrinsic: state stored inside the flyweight, is independent of the context in which its used."
